$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M2").Value = -1.833336000000003
$ws.Range("H2").Value = 114.833336
$ws.Range("K2").Value = 114.833336
$ws.Range("I2").Value = 114.833336
$ws.Range("I6").Value = 2524.75
$ws.Range("K6").Value = 7574.25
$ws.Range("M6").Value = -7462.25
$ws.Range("H6").Value = 3946.4666
$ws.Range("N9").ClearContents()
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -11082.111
$ws.Range("H9").Value = 11251.111
$ws.Range("I9").Value = 11251.111
$ws.Range("K9").Value = 11251.111
$ws.Range("J9").Value = 0
$ws.Range("H19").Value = 1065.7693
$ws.Range("I19").Value = 1033.125
$ws.Range("K19").Value = 1033.125
$ws.Range("M19").Value = -858.125
$ws.Range("N19").Value = -1468
$ws.Range("J19").Value = 1118
$ws.Range("L19").Value = 1118
$ws.Range("K40").Value = 2137.4
$ws.Range("I40").Value = 2137.4
$ws.Range("H40").Value = 3126.6428
$ws.Range("M40").Value = -1962.4
$ws.Range("L112").Value = 609595.2
$ws.Range("N112").Value = -611811.2
$ws.Range("J112").Value = 203198.4
$ws.Range("H112").Value = 201979
$ws.Range("M137").Value = -936.6000000000004
$ws.Range("K137").Value = 3486.6
$ws.Range("J137").Value = 2164.889
$ws.Range("H137").Value = 1538.2084
$ws.Range("L137").Value = 6494.667
$ws.Range("N137").Value = -11594.667
$ws.Range("I137").Value = 1162.2
$ws.Range("N138").Value = -46420.25
$ws.Range("H138").Value = 5471.8
$ws.Range("J138").Value = 12046.75
$ws.Range("L138").Value = 36140.25
$ws.Range("I141").Value = 3539.75
$ws.Range("K141").Value = 10619.25
$ws.Range("M141").Value = -5439.25
$ws.Range("H141").Value = 4256.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -417
$ws.Range("H2").Value = 1743.1
$ws.Range("K2").Value = 530
$ws.Range("I2").Value = 530
$ws.Range("L5").Value = 3520
$ws.Range("M5").Value = -506.6667
$ws.Range("K5").Value = 618.6667
$ws.Range("H5").Value = 2496
$ws.Range("I5").Value = 618.6667
$ws.Range("N5").Value = -3744
$ws.Range("J5").Value = 3520
$ws.Range("I32").Value = 1067.8452
$ws.Range("M32").Value = -780.8452
$ws.Range("H32").Value = 1825.2307
$ws.Range("K32").Value = 1067.8452
$ws.Range("J45").Value = 2778
$ws.Range("I45").Value = 8916.333000000001
$ws.Range("M45").Value = -8539.333000000001
$ws.Range("H45").Value = 6461
$ws.Range("K45").Value = 8916.333000000001
$ws.Range("L45").Value = 2778
$ws.Range("N45").Value = -3532
$ws.Range("I61").Value = 37042252
$ws.Range("K61").Value = 37042252
$ws.Range("N61").Value = -7423.5
$ws.Range("L61").Value = 6999.5
$ws.Range("M61").Value = -37042040
$ws.Range("H61").Value = 34488096
$ws.Range("J61").Value = 6999.5
$ws.Range("I97").Value = 437.6
$ws.Range("M97").Value = 58.39999999999998
$ws.Range("H97").Value = 437.6
$ws.Range("K97").Value = 437.6
$ws.Range("H110").Value = 73315.92999999999
$ws.Range("K110").Value = 101041.1
$ws.Range("M110").Value = -98996.10000000001
$ws.Range("I110").Value = 101041.1
$ws.Range("M116").Value = 1764
$ws.Range("I116").Value = 530
$ws.Range("K116").Value = 530
$ws.Range("H116").Value = 1743.1
$ws.Range("H122").Value = 7168.4287
$ws.Range("K122").Value = 17964.4995
$ws.Range("I122").Value = 5988.1665
$ws.Range("M122").Value = -15514.4995
$ws.Range("N136").Value = -26098.5
$ws.Range("M136").Value = -111124206
$ws.Range("L136").Value = 20998.5
$ws.Range("H136").Value = 34488096
$ws.Range("K136").Value = 111126756
$ws.Range("I136").Value = 37042252
$ws.Range("J136").Value = 6999.5
$ws.Range("J141").Value = 92645
$ws.Range("L141").Value = 92645
$ws.Range("N141").Value = -103005
$ws.Range("H141").Value = 92645

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -416
$ws.Range("I3").Value = 530
$ws.Range("H3").Value = 1743.1
$ws.Range("K3").Value = 530
$ws.Range("H4").Value = 2496
$ws.Range("I4").Value = 618.6667
$ws.Range("K4").Value = 618.6667
$ws.Range("N4").Value = -3750
$ws.Range("L4").Value = 3520
$ws.Range("J4").Value = 3520
$ws.Range("M4").Value = -503.6667
$ws.Range("H94").Value = 1723.85
$ws.Range("M94").Value = -1258.3158
$ws.Range("I94").Value = 1709.3158
$ws.Range("N94").Value = -2902
$ws.Range("K94").Value = 1709.3158
$ws.Range("J94").Value = 2000
$ws.Range("L94").Value = 2000
$ws.Range("M99").Value = -1201.3333
$ws.Range("H99").Value = 2725.75
$ws.Range("K99").Value = 2699.3333
$ws.Range("I99").Value = 2699.3333
$ws.Range("H107").Value = 112206.664
$ws.Range("M107").Value = 1264.63635
$ws.Range("K107").Value = 655.36365
$ws.Range("I107").Value = 655.36365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L7").Value = 224.14285
$ws.Range("K7").Value = 280
$ws.Range("N7").Value = -450.14285
$ws.Range("H7").Value = 253.93333
$ws.Range("I7").Value = 280
$ws.Range("J7").Value = 224.14285
$ws.Range("M7").Value = -167
$ws.Range("I31").Value = 12712.333
$ws.Range("M31").Value = -12417.333
$ws.Range("J31").Value = 2656.8333
$ws.Range("L31").Value = 2656.8333
$ws.Range("H31").Value = 7684.5835
$ws.Range("K31").Value = 12712.333
$ws.Range("N31").Value = -3246.8333
$ws.Range("I34").Value = 12712.333
$ws.Range("H34").Value = 7684.5835
$ws.Range("N34").Value = -3060.8333
$ws.Range("J34").Value = 2656.8333
$ws.Range("M34").Value = -12510.333
$ws.Range("K34").Value = 12712.333
$ws.Range("L34").Value = 2656.8333
$ws.Range("H58").Value = 17864276
$ws.Range("I58").Value = 23817764
$ws.Range("K58").Value = 23817764
$ws.Range("M58").Value = -23817561
$ws.Range("M136").Value = -71450742
$ws.Range("H136").Value = 17864276
$ws.Range("K136").Value = 71453292
$ws.Range("I136").Value = 23817764
$ws.Range("N140").Value = -77329.66
$ws.Range("J140").Value = 66969.66
$ws.Range("L140").Value = 66969.66
$ws.Range("H140").Value = 66969.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 19799.8
$ws.Range("I76").Value = 18999
$ws.Range("M76").Value = -56614
$ws.Range("K76").Value = 56997
$ws.Range("K79").Value = 56997
$ws.Range("M79").Value = -55671
$ws.Range("H79").Value = 19799.8
$ws.Range("I79").Value = 18999
$ws.Range("K121").Value = 232029.702
$ws.Range("H121").Value = 48792.56
$ws.Range("J121").Value = 17862.666
$ws.Range("L121").Value = 53587.99800000001
$ws.Range("I121").Value = 77343.234
$ws.Range("M121").Value = -230719.702
$ws.Range("N121").Value = -56207.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L17").Value = 17504
$ws.Range("J17").Value = 17504
$ws.Range("H17").Value = 17504
$ws.Range("N17").Value = -17840
$ws.Range("N43").ClearContents()
$ws.Range("L43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("K57").Value = 4966.6665
$ws.Range("H57").Value = 4966.6665
$ws.Range("M57").Value = -4146.6665
$ws.Range("J57").Value = 0
$ws.Range("I57").Value = 4966.6665
$ws.Range("L57").Value = 0
$ws.Range("K80").Value = 1999
$ws.Range("M80").Value = -1001
$ws.Range("J80").Value = 1899.5
$ws.Range("L80").Value = 1899.5
$ws.Range("N80").Value = -3895.5
$ws.Range("I80").Value = 1999
$ws.Range("H80").Value = 1949.25
$ws.Range("J83").Value = 1899.5
$ws.Range("N83").Value = -19481.5
$ws.Range("M83").Value = -5003
$ws.Range("K83").Value = 9995
$ws.Range("H83").Value = 1949.25
$ws.Range("I83").Value = 1999
$ws.Range("L83").Value = 9497.5
$ws.Range("I113").Value = 59566.945
$ws.Range("H113").Value = 53860.05
$ws.Range("K113").Value = 59566.945
$ws.Range("M113").Value = -57396.945
$ws.Range("H117").Value = 65000
$ws.Range("N117").Value = -71884
$ws.Range("L117").Value = 65000
$ws.Range("J117").Value = 65000
$ws.Range("M132").Value = -16305773
$ws.Range("H132").Value = 4630887.5
$ws.Range("K132").Value = 16308303
$ws.Range("I132").Value = 5436101

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2379.5454
$ws.Range("N22").Value = -2377.4
$ws.Range("K22").Value = 2873
$ws.Range("L22").Value = 1787.4
$ws.Range("I22").Value = 2873
$ws.Range("J22").Value = 1787.4
$ws.Range("M22").Value = -2578
$ws.Range("I27").Value = 2873
$ws.Range("M27").Value = -2766
$ws.Range("N27").Value = -2001.4
$ws.Range("H27").Value = 2379.5454
$ws.Range("L27").Value = 1787.4
$ws.Range("K27").Value = 2873
$ws.Range("J27").Value = 1787.4
$ws.Range("L46").Value = 795.6667
$ws.Range("J46").Value = 795.6667
$ws.Range("H46").Value = 1143.7273
$ws.Range("N46").Value = -1171.6667
$ws.Range("I61").Value = 3468
$ws.Range("K61").Value = 3468
$ws.Range("N61").Value = -2404
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -3266
$ws.Range("H61").Value = 3390.7368
$ws.Range("J61").Value = 2000
$ws.Range("I113").Value = 3468
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("H113").Value = 3390.7368
$ws.Range("K113").Value = 3468
$ws.Range("N113").Value = -6340
$ws.Range("M113").Value = -1298
$ws.Range("N136").Value = -27595.0005
$ws.Range("M136").Value = -14682.75
$ws.Range("L136").Value = 22495.0005
$ws.Range("H136").Value = 6496
$ws.Range("K136").Value = 17232.75
$ws.Range("I136").Value = 5744.25
$ws.Range("J136").Value = 7498.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("J109").Value = 0
$ws.Range("H109").Value = 0
$ws.Range("M132").Value = -100001384
$ws.Range("H132").Value = 29413098
$ws.Range("K132").Value = 100003914
$ws.Range("I132").Value = 33334638
